$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 41
$ws.Range("F3").Value = 30
$ws.Range("H3").Value = 31

$ws.Range("E4").Value = 45
$ws.Range("F4").Value = 28
$ws.Range("H4").Value = 40

$ws.Range("E6").Value = 44

$ws.Range("E8").Value = 8

$ws.Range("E10").Value = 568
$ws.Range("F10").Value = 280
$ws.Range("G10").Value = 96
$ws.Range("H10").Value = 376

$ws.Range("E11").Value = 361
$ws.Range("G11").Value = 64
$ws.Range("H11").Value = 259

$ws.Range("E12").Value = 552
$ws.Range("F12").Value = 295
$ws.Range("H12").Value = 381

$ws.Range("E13").Value = 139

$ws.Range("E15").Value = 168
$ws.Range("F15").Value = 73
$ws.Range("H15").Value = 123

$ws.Range("E16").Value = 204
$ws.Range("F16").Value = 104
$ws.Range("H16").Value = 152

$ws.Range("E19").Value = 14

$ws.Range("E20").Value = 87

$ws.Range("E21").Value = 139

$ws.Range("E22").Value = 173
$ws.Range("F22").Value = 92
$ws.Range("H22").Value = 134

$ws.Range("E23").Value = 203

$ws.Range("E25").Value = 269

$ws.Range("E27").Value = 333
$ws.Range("F27").Value = 170
$ws.Range("H27").Value = 251

$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 81
$ws.Range("H28").Value = 133

$ws.Range("F29").Value = 95
$ws.Range("H29").Value = 136

$ws.Range("E32").Value = 185
$ws.Range("F32").Value = 109
$ws.Range("H32").Value = 147

$ws.Range("E33").Value = 294
$ws.Range("F33").Value = 151
$ws.Range("H33").Value = 240

$ws.Range("E34").Value = 220
$ws.Range("F34").Value = 146
$ws.Range("H34").Value = 185

$ws.Range("E35").Value = 152
$ws.Range("F35").Value = 97
$ws.Range("H35").Value = 124

$ws.Range("E36").Value = 74
$ws.Range("F36").Value = 43
$ws.Range("H36").Value = 53

$ws.Range("E37").Value = 164
$ws.Range("F37").Value = 78
$ws.Range("H37").Value = 115

$ws.Range("E40").Value = 262

$ws.Range("E41").Value = 393
$ws.Range("F41").Value = 186
$ws.Range("H41").Value = 278

$ws.Range("E42").Value = 383
$ws.Range("F42").Value = 215
$ws.Range("H42").Value = 276

$ws.Range("E44").Value = 313

$ws.Range("E45").Value = 146

$ws.Range("E46").Value = 325

$ws.Range("E47").Value = 460
$ws.Range("F47").Value = 234
$ws.Range("H47").Value = 326

$ws.Range("E48").Value = 215

$ws.Range("E49").Value = 286
$ws.Range("F49").Value = 126
$ws.Range("H49").Value = 213

$ws.Range("E50").Value = 244

$ws.Range("E51").Value = 238
$ws.Range("F51").Value = 106
$ws.Range("H51").Value = 179

$ws.Range("E52").Value = 28
